$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.206.10"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -2.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.873.39"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.67%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.90"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.71%  "
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5043"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3752"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07164"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8898"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.77"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.876.32"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07579"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.334"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.57"
$ws.Range("D15").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008540"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.16"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.256.70"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.085"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.114.17"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.64"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.48%  "
$ws.Range("E24").Value = "  -1.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.06"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.840"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.03"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.83%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.085"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -5.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.81"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.763"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.713"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08981"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05144"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.106"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7464"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.163"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.552"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02030"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.038"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("E40").Value = "  -1.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5368"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.638"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "114.87"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.451"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1479"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4648"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.003"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("E48").Value = "  -4.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.571"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "64.88"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.87"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.32%  "
